$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 53 (ALC)
$ws.Range("H53").Value = 229.125
$ws.Range("I53").Value = 94.875
$ws.Range("J53").Value = 363.375
$ws.Range("K53").Value = 94.875
$ws.Range("L53").Value = 363.375
$ws.Range("M53").Value = 542.125
$ws.Range("N53").Value = -1637.375

# Row 58 (ALC)
$ws.Range("H58").Value = 638.3333
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

# Row 107 (ALC)
$ws.Range("H107").Value = 745.4286
$ws.Range("I107").Value = 637.2222
$ws.Range("K107").Value = 637.2222
$ws.Range("M107").Value = 1282.7778

# Row 132 (ALC)
$ws.Range("H132").Value = 3904.158
$ws.Range("I132").Value = 2477.9
$ws.Range("K132").Value = 7433.700000000001
$ws.Range("M132").Value = -4903.700000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (ARM)
$ws.Range("H61").Value = 2296.9
$ws.Range("J61").Value = 2322
$ws.Range("L61").Value = 2322
$ws.Range("N61").Value = -2746

# Row 109 (ARM)
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774

# Row 122 (ARM)
$ws.Range("H122").Value = 1935.5
$ws.Range("I122").Value = 1943.6
$ws.Range("K122").Value = 5830.799999999999
$ws.Range("M122").Value = -3380.799999999999

# Row 132 (ARM)
$ws.Range("H132").Value = 6398
$ws.Range("I132").Value = 2499.5
$ws.Range("J132").Value = 8997
$ws.Range("K132").Value = 7498.5
$ws.Range("L132").Value = 26991
$ws.Range("M132").Value = -4968.5
$ws.Range("N132").Value = -32051

# Row 136 (ARM)
$ws.Range("H136").Value = 2296.9
$ws.Range("J136").Value = 2322
$ws.Range("L136").Value = 6966
$ws.Range("N136").Value = -12066

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (BSM)
$ws.Range("H86").Value = 2892.4443
$ws.Range("I86").Value = 2254.1667
$ws.Range("J86").Value = 4169
$ws.Range("K86").Value = 2254.1667
$ws.Range("L86").Value = 4169
$ws.Range("M86").Value = -1131.1667
$ws.Range("N86").Value = -6415

# Row 89 (BSM)
$ws.Range("H89").Value = 2892.4443
$ws.Range("I89").Value = 2254.1667
$ws.Range("J89").Value = 4169
$ws.Range("K89").Value = 11270.8335
$ws.Range("L89").Value = 20845
$ws.Range("M89").Value = -5654.833500000001
$ws.Range("N89").Value = -32077

$ws = $wb.Worksheets.Item("CRP")
# Row 58 (CRP)
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

# Row 86 (CRP)
$ws.Range("H86").Value = 8743
$ws.Range("I86").Value = 7995
$ws.Range("K86").Value = 7995
$ws.Range("M86").Value = -6872

# Row 89 (CRP)
$ws.Range("H89").Value = 8743
$ws.Range("I89").Value = 7995
$ws.Range("K89").Value = 39975
$ws.Range("M89").Value = -34359

# Row 99 (CRP)
$ws.Range("H99").Value = 4734.4346
$ws.Range("I99").Value = 4345.8823
$ws.Range("J99").Value = 5835.3335
$ws.Range("K99").Value = 4345.8823
$ws.Range("L99").Value = 5835.3335
$ws.Range("M99").Value = -2847.8823
$ws.Range("N99").Value = -8831.333500000001

# Row 126 (CRP)
$ws.Range("H126").Value = 4734.4346
$ws.Range("I126").Value = 4345.8823
$ws.Range("J126").Value = 5835.3335
$ws.Range("K126").Value = 13037.6469
$ws.Range("L126").Value = 17506.0005
$ws.Range("M126").Value = -10567.6469
$ws.Range("N126").Value = -22446.0005

# Row 134 (CRP)
$ws.Range("H134").Value = 3921.7778
$ws.Range("I134").Value = 3412.25
$ws.Range("K134").Value = 10236.75
$ws.Range("M134").Value = -7701.75

# Row 136 (CRP)
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 34 (CUL)
$ws.Range("H34").Value = 484.375
$ws.Range("J34").Value = 1008.3333
$ws.Range("L34").Value = 3024.9999
$ws.Range("N34").Value = -3192.9999

# Row 39 (CUL)
$ws.Range("H39").Value = 2771.8572
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 2771.8572
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 8315.571599999999
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -8903.571599999999

# Row 55 (CUL)
$ws.Range("H55").Value = 1715.8889
$ws.Range("I55").Value = 497.66666
$ws.Range("J55").Value = 2325
$ws.Range("K55").Value = 1492.99998
$ws.Range("L55").Value = 6975
$ws.Range("M55").Value = -1315.99998
$ws.Range("N55").Value = -7329

# Row 64 (CUL)
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

# Row 67 (CUL)
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 15 (GSM)
$ws.Range("H15").Value = 22996.111
$ws.Range("J15").Value = 22996.111
$ws.Range("L15").Value = 22996.111
$ws.Range("N15").Value = -23572.111

# Row 43 (GSM)
$ws.Range("H43").Value = 8599.666999999999
$ws.Range("J43").Value = 20000
$ws.Range("L43").Value = 20000
$ws.Range("N43").Value = -20302

# Row 81 (GSM)
$ws.Range("H81").Value = 22996.111
$ws.Range("J81").Value = 22996.111
$ws.Range("L81").Value = 22996.111
$ws.Range("N81").Value = -24992.111

# Row 84 (GSM)
$ws.Range("H84").Value = 22996.111
$ws.Range("J84").Value = 22996.111
$ws.Range("L84").Value = 68988.333
$ws.Range("N84").Value = -78972.333

# Row 102 (GSM)
$ws.Range("H102").Value = 445.07144
$ws.Range("I102").Value = 460.9
$ws.Range("J102").Value = 405.5
$ws.Range("K102").Value = 460.9
$ws.Range("L102").Value = 405.5
$ws.Range("M102").Value = 1161.1
$ws.Range("N102").Value = -3649.5

# Row 123 (GSM)
$ws.Range("H123").Value = 35666
$ws.Range("J123").Value = 35666
$ws.Range("L123").Value = 35666
$ws.Range("N123").Value = -40566

$ws = $wb.Worksheets.Item("LTW")
# Row 61 (LTW)
$ws.Range("H61").Value = 1759.7858
$ws.Range("I61").Value = 1420.1818
$ws.Range("K61").Value = 1420.1818
$ws.Range("M61").Value = -1218.1818

# Row 93 (LTW)
$ws.Range("H93").Value = 975.375
$ws.Range("I93").Value = 900.4286
$ws.Range("J93").Value = 1500
$ws.Range("K93").Value = 900.4286
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = 347.5714
$ws.Range("N93").Value = -3996

# Row 113 (LTW)
$ws.Range("H113").Value = 1759.7858
$ws.Range("I113").Value = 1420.1818
$ws.Range("K113").Value = 1420.1818
$ws.Range("M113").Value = 749.8181999999999

# Row 136 (LTW)
$ws.Range("H136").Value = 4581.7856
$ws.Range("I136").Value = 4477.273
$ws.Range("J136").Value = 4965
$ws.Range("K136").Value = 13431.819
$ws.Range("L136").Value = 14895
$ws.Range("M136").Value = -10881.819
$ws.Range("N136").Value = -19995

$ws = $wb.Worksheets.Item("WVR")
# Row 75 (WVR)
$ws.Range("H75").Value = 24999.5

# Row 78 (WVR)
$ws.Range("H78").Value = 24999.5

# Row 107 (WVR)
$ws.Range("H107").Value = 700.0769
$ws.Range("I107").Value = 719.1667
$ws.Range("J107").Value = 683.7143
$ws.Range("K107").Value = 2157.5001
$ws.Range("L107").Value = 2051.1429
$ws.Range("M107").Value = -237.5001000000002
$ws.Range("N107").Value = -5891.1429

# Row 113 (WVR)
$ws.Range("H113").Value = 621.2857
$ws.Range("I113").Value = 392.75
$ws.Range("J113").Value = 926
$ws.Range("K113").Value = 1178.25
$ws.Range("L113").Value = 2778
$ws.Range("M113").Value = 991.75
$ws.Range("N113").Value = -7118

# Row 122 (WVR)
$ws.Range("H122").Value = 2142.889
$ws.Range("I122").Value = 2142.889
$ws.Range("K122").Value = 6428.667
$ws.Range("M122").Value = -3978.667
